$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G12").Value = 1240524717.0500007
$ws.Range("I12").Value = 4188377156

$ws.Range("G13").Value = 319819483.18000001
$ws.Range("I13").Value = 1012006300

$ws.Range("G14").Value = 34063116.800000042
$ws.Range("I14").Value = 44319159.289999999

$ws.Range("G16").Value = -60834434.380000003
$ws.Range("I16").Value = 162861896.59999999

$ws.Range("G18").Formula = "=SUM(G12:G17)"

$ws.Range("G19").Value = -379300000.00000012
$ws.Range("I19").Value = 1160500000

$ws.Range("G21").Formula = "=SUM(G18:G20)"

$ws.Range("I22").Value = 85592745

$ws.Range("G26").Value = 1029174575.116062
$ws.Range("I26").Value = 1010658959

$wb.Application.CalculateFull()
